$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a brand-new column A (everything else shifts right by one) ---
$ws.Columns("A:A").Insert()

# --- Capture the old values that need to move to new homes before we ---
# --- overwrite the cells they currently sit in (columns already      ---
# --- shifted right by one from the Insert() above). Old D2's date    ---
# --- text is dropped in favor of a real Date value written into B2.  ---
$oldLocation = $ws.Range("B2").Value2   # was old A2 ("vault_lake")
$oldStart    = $ws.Range("C2").Value2   # was old B2 (start time fraction)
$oldStop     = $ws.Range("D2").Value2   # was old C2 (stop time fraction)

# --- New "program run?" column (A): header + its values first ---
$ws.Range("A1").Value2 = "program_run?"
$ws.Range("A2").Value2 = "y"
$ws.Range("A4").Value2 = "p"
$ws.Range("A5").Value2 = "y"

# --- Remaining new headers for columns B-G ---
$ws.Range("B1").Value2 = "date_(yyyy-mm-dd)"
$ws.Range("C1").Value2 = "location_(lake)"
$ws.Range("D1").Value2 = "start_time_(hh:mm:ss)"
$ws.Range("E1").Value2 = "stop_time_(hh:mm:ss)"
$ws.Range("F1").Value2 = "surface_type"
$ws.Range("G1").Value2 = "surface_class"

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 10
$headerRange.Font.Name = "Arial"
$headerRange.WrapText = $true

# The old "Program run?" header used to live at the end (old AW1); after
# the column insert it landed at AX1 - that text now lives at A1 instead,
# so drop the now-empty trailing column entirely.
$ws.Columns("AX:AX").Delete()

# --- Re-point the row-2 sample data into its new column layout ---
$ws.Range("C2").ClearFormats()
$ws.Range("C2").Value2 = $oldLocation
$ws.Range("D2").Value2 = $oldStart
$ws.Range("D2").NumberFormat = "h:mm:ss"
$ws.Range("E2").Value2 = $oldStop
$ws.Range("E2").NumberFormat = "h:mm:ss"

$ws.Range("B2").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B2").Value2 = 44042

# --- New rows of data below ---
$ws.Range("B3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B3").Value2 = 44042

$ws.Range("B4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B4").Value2 = 44042

$ws.Range("B5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B5").Value2 = 44042

$ws.Range("B6").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B6").Value2 = 44042

# --- Column widths: B keeps the old "date" column width, C:D keep the ---
# --- old "location/start-time" width, the rest just shift down one.  ---
$ws.Columns("B:B").ColumnWidth = 12.608072916666666
$ws.Columns("C:D").ColumnWidth = 14.608072916666666
$ws.Columns("E:E").ColumnWidth = 21.05
$ws.Columns("F:J").ColumnWidth = 11.83
$ws.Columns("K:K").ColumnWidth = 24.608072916666666
$ws.Columns("L:L").ColumnWidth = 22.16
$ws.Columns("M:M").ColumnWidth = 12.83
$ws.Columns("O:O").ColumnWidth = 15.83
$ws.Columns("P:P").ColumnWidth = 12.05
$ws.Columns("R:R").ColumnWidth = 14.72
$ws.Columns("AT:AT").ColumnWidth = 12.05
$ws.Columns("AU:AU").ColumnWidth = 11.83
$ws.Columns("AV:AV").ColumnWidth = 13.17
$ws.Columns("AW:AW").ColumnWidth = 11.94

# --- active selection moved in the authored file ---
$ws.Range("G8").Select()
